# "finished analysing the data"
#
# The "Record of Progress" table gets a new entry for 07/10/2017. The
# previous last row (05/10/2017, "I analysed more data") keeps its text
# but loses the hidden _GoBack bookmark (which Word always keeps pinned
# to the most recent edit) and gains a trailing blank paragraph in its
# "Task" cell. The new row carries the "finished analysing the data" /
# "will need to go back and rephrase..." text, and the _GoBack bookmark
# now lives at the end of its last cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Step 1: Strip the old _GoBack bookmark out of the current last row.
# Deleting the paragraph's range (text + the invisible bookmark anchor
# riding along at its end) and retyping the text is the reliable way to
# drop the bookmark without leaving an orphaned copy behind.
# ---------------------------------------------------------------------
$lastRow = $t.Rows.Item($t.Rows.Count)
$taskCell = $lastRow.Cells.Item(2)
$taskCell.Range.Paragraphs.Item(1).Range.Delete()
$taskCell.Range.Text = "I analysed more data"
$taskCell.Range.Paragraphs.Item($taskCell.Range.Paragraphs.Count).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# Step 2: Add the new progress-log row for 07/10/2017.
# ---------------------------------------------------------------------
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "07/10/2017"
$newRow.Cells.Item(2).Range.Text = "I finished analysing the data"

$commentsCell = $newRow.Cells.Item(3)
$commentsCell.Range.Text = "I will need to go back and rephrase some of the analyses and conclusions"
$commentsCell.Range.Paragraphs.Item($commentsCell.Range.Paragraphs.Count).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# Step 3: Re-plant the _GoBack bookmark at the end of the new last cell,
# matching where Word leaves it after the most recent edit.
# ---------------------------------------------------------------------
$newLastPara = $commentsCell.Range.Paragraphs.Item($commentsCell.Range.Paragraphs.Count)
$goBackRange = $newLastPara.Range
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

Write-Output "Table now has $($t.Rows.Count) rows."
